$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old "Category" column D -> E,
# and the D:D merge for "Category" becomes E2:E4).
$ws.Columns("D:D").Insert()

# New column D, row 3 gets a "Heart" header, merged over D3:D4 (mirrors the
# existing B4/C4 "Age"/"Birth" sub-row pattern under the B3:C3 "Generation" cell).
$ws.Range("D3").Value = "Heart"
$ws.Range("D3:D4").Merge()

# The "Life" header (originally merged B2:C2) now spans the new column too.
$ws.Range("B2:D2").Merge()
